# Updated cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text (e.g. "27.176.70"),
# so force text formatting before writing to avoid Excel auto-converting
# numeric-looking strings into numbers; formatting is cleared afterwards
# so the cells keep their original (unstyled) appearance.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.176.70'
$ws.Range("E2").Value = '  -0.42%  '

# Row 3
$ws.Range("D3").Value = '1.700.60'
$ws.Range("E3").Value = '  -0.45%  '

# Row 4
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.55%  '

# Row 5
$ws.Range("D5").Value = '223.57'
$ws.Range("E5").Value = '  -0.02%  '

# Row 6
$ws.Range("D6").Value = '0.5242'
$ws.Range("E6").Value = '  -0.92%  '

# Row 7
$ws.Range("E7").Value = '  +0.58%  '

# Row 8
$ws.Range("D8").Value = '0.06599'
$ws.Range("E8").Value = '  +1.00%  '

# Row 9
$ws.Range("D9").Value = '0.2618'
$ws.Range("E9").Value = '  -0.76%  '

# Row 10
$ws.Range("D10").Value = '20.54'
$ws.Range("E10").Value = '  -1.74%  '

# Row 11
$ws.Range("D11").Value = '0.07706'
$ws.Range("E11").Value = '  +1.00%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.710.16'
$ws.Range("E12").Value = '  +0.03%  '

# Row 13
$ws.Range("D13").Value = '1.936.64'
$ws.Range("E13").Value = '  -0.39%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.408'
$ws.Range("E14").Value = '  -3.43%  '

# Row 15
$ws.Range("D15").Value = '0.5728'
$ws.Range("E15").Value = '  -0.20%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8111'
$ws.Range("E16").Value = '  -0.96%  '

# Row 17
$ws.Range("D17").Value = '66.81'
$ws.Range("E17").Value = '  -0.62%  '

# Row 18
$ws.Range("D18").Value = '27.206.22'
$ws.Range("E18").Value = '  -0.29%  '

# Row 19
$ws.Range("D19").Value = '217.35'
$ws.Range("E19").Value = '  +0.94%  '

# Row 20
$ws.Range("D20").Value = '1.010'
$ws.Range("E20").Value = '  +0.70%  '

# Row 21
$ws.Range("D21").Value = '4.581'
$ws.Range("E21").Value = '  -1.81%  '

# Row 22
$ws.Range("E22").Value = '  -1.37%  '

# Row 23
$ws.Range("D23").Value = '5.995'
$ws.Range("E23").Value = '  +0.56%  '

# Row 24
$ws.Range("D24").Value = '1.010'
$ws.Range("E24").Value = '  +0.58%  '

# Row 25
$ws.Range("D25").Value = '145.12'
$ws.Range("E25").Value = '  +1.94%  '

# Row 26
$ws.Range("D26").Value = '1.740'
$ws.Range("E26").Value = '  -0.72%  '

# Row 27
$ws.Range("D27").Value = '0.1195'
$ws.Range("E27").Value = '  -1.78%  '

# Row 28
$ws.Range("D28").Value = '7.147'
$ws.Range("E28").Value = '  -1.59%  '

# Row 29
$ws.Range("D29").Value = '16.01'
$ws.Range("E29").Value = '  -1.76%  '

# Row 30
$ws.Range("D30").Value = '0.05288'
$ws.Range("E30").Value = '  -1.42%  '

# Row 31
$ws.Range("D31").Value = '1.287'
$ws.Range("E31").Value = '  -0.41%  '

# Row 32
$ws.Range("D32").Value = '3.419'
$ws.Range("E32").Value = '  -1.89%  '

# Row 33
$ws.Range("D33").Value = '3.309'
$ws.Range("E33").Value = '  -3.11%  '

# Row 34
$ws.Range("D34").Value = '1.622'
$ws.Range("E34").Value = '  -1.02%  '

# Row 35
$ws.Range("D35").Value = '2.814'
$ws.Range("E35").Value = '  -2.00%  '

# Row 36
$ws.Range("D36").Value = '2.402'
$ws.Range("E36").Value = '  -0.64%  '

# Row 37
$ws.Range("D37").Value = '0.9391'
$ws.Range("E37").Value = '  -1.10%  '

# Row 38
$ws.Range("D38").Value = '0.5819'
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
$ws.Range("D39").Value = '1.180.95'
$ws.Range("E39").Value = '  +13.39%  '

# Row 40
$ws.Range("E40").Value = '  +0.52%  '

# Row 41
$ws.Range("D41").Value = '1.009'
$ws.Range("E41").Value = '  +0.61%  '

# Row 42
$ws.Range("D42").Value = '5.725'
$ws.Range("E42").Value = '  -2.46%  '

# Row 43
$ws.Range("D43").Value = '0.8354'
$ws.Range("E43").Value = '  -0.44%  '

# Row 44
$ws.Range("D44").Value = '100.81'
$ws.Range("E44").Value = '  -0.21%  '

# Row 45
$ws.Range("D45").Value = '1.849.01'
$ws.Range("E45").Value = '  -0.18%  '

# Row 46
$ws.Range("D46").Value = '0.0₈109'
$ws.Range("E46").Value = '  -5.05%  '

# Row 47
$ws.Range("D47").Value = '56.98'
$ws.Range("E47").Value = '  -1.79%  '

# Row 48
$ws.Range("D48").Value = '0.4554'
$ws.Range("E48").Value = '  +1.30%  '

# Row 49
$ws.Range("D49").Value = '1.008'
$ws.Range("E49").Value = '  +0.53%  '

# Row 50
$ws.Range("D50").Value = '8.056'
$ws.Range("E50").Value = '  -0.36%  '

# Row 51
$ws.Range("D51").Value = '0.05220'
$ws.Range("E51").Value = '  -0.31%  '

$priceRange.ClearFormats()
